$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "Modelo", matching the formatting of the existing
# header row (E1) by copying its format.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F1").Value = "Modelo"

# Update the numeric prediction values in row 2
$ws.Range("B2").Value = 0.4161312050481046
$ws.Range("C2").Value = 0.9917944623708683
$ws.Range("D2").Value = 0.4987128619250442

# Add the model name in the new F2 cell
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=3, n_estimators=50))])"
